$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @{ Row = 2; B = "Bitcoin"; C = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"; D = "63.067.81"; E = "  -2.44%  " },
  @{ Row = 3; B = "Ethereum"; C = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"; D = "3.220.74"; E = "  -4.06%  " },
  @{ Row = 4; B = "TetherUSD"; C = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"; D = "0.999"; E = "  -0.07%  " },
  @{ Row = 5; B = "BNB"; C = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"; D = "533.80"; E = "  +1.06%  " },
  @{ Row = 6; B = "Solana"; C = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"; D = "173.97"; E = "  -6.41%  " },
  @{ Row = 7; B = "XRP"; C = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"; D = "0.596"; E = "  -0.29%  " },
  @{ Row = 8; B = "USDC"; C = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"; D = "1.00"; E = "  -0.02%  " },
  @{ Row = 9; B = "LidoStakedEther"; C = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"; D = "3.220.23"; E = "  -3.93%  " },
  @{ Row = 10; B = "Cardano"; C = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"; D = "0.611"; E = "  -1.72%  " },
  @{ Row = 11; B = "Avalanche"; C = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; D = "54.05"; E = "  -9.20%  " },
  @{ Row = 12; B = "Dogecoin"; C = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"; D = "0.134"; E = "  +1.21%  " },
  @{ Row = 13; B = "ShibaInu"; C = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; D = "0.0000254"; E = "  +0.02%  " },
  @{ Row = 14; B = "Polkadot"; C = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"; D = "9.13"; E = "  -0.44%  " },
  @{ Row = 15; B = "WrappedliquidstakedEther2.0"; C = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; D = "3.723.00"; E = "  -4.24%  " },
  @{ Row = 16; B = "TRON"; C = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; D = "0.117"; E = "  -3.29%  " },
  @{ Row = 17; B = "Chainlink"; C = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; D = "17.41"; E = "  +0.69%  " },
  @{ Row = 18; B = "WrappedEther"; C = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D = "3.208.29"; E = "  -4.13%  " },
  @{ Row = 19; B = "WrappedBTC"; C = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D = "62.780.06"; E = "  -2.49%  " },
  @{ Row = 20; B = "Uniswap"; C = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; D = "11.12"; E = "  +1.03%  " },
  @{ Row = 21; B = "Polygon"; C = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"; D = "0.971"; E = "  +1.43%  " },
  @{ Row = 22; B = "BitcoinCash"; C = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D = "368.54"; E = "  -1.63%  " },
  @{ Row = 23; B = "PancakeSwap"; C = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"; D = "3.79"; E = "  +1.71%  " },
  @{ Row = 24; B = "RenderToken"; C = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; D = "11.31"; E = "  +4.04%  " },
  @{ Row = 25; B = "Litecoin"; C = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; D = "81.48"; E = "  +1.25%  " },
  @{ Row = 26; B = "Toncoin"; C = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; D = "3.89"; E = "  +2.48%  " },
  @{ Row = 27; B = "LEO"; C = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"; D = "6.13"; E = "  +1.71%  " },
  @{ Row = 28; B = "ImmutableX"; C = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D = "2.68"; E = "  +1.15%  " },
  @{ Row = 29; B = "InternetComputer(DFINITY)"; C = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D = "11.39"; E = "  +0.67%  " },
  @{ Row = 30; B = "Filecoin"; C = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D = "8.23"; E = "  -1.96%  " },
  @{ Row = 31; B = "EthereumClassic"; C = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D = "28.59"; E = "  -0.98%  " },
  @{ Row = 32; B = "Bittensor"; C = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"; D = "646.57"; E = "  -1.22%  " },
  @{ Row = 33; B = "NEARProtocol"; C = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; D = "6.58"; E = "  -3.02%  " },
  @{ Row = 34; B = "Cosmos"; C = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; D = "11.39"; E = "  +2.36%  " },
  @{ Row = 35; B = "Hedera"; C = "https://coinranking.com/coin/jad286TjB+hedera-hbar"; D = "0.106"; E = "  +2.28%  " },
  @{ Row = 36; B = "OKB"; C = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"; D = "57.20"; E = "  -4.19%  " },
  @{ Row = 37; B = "Dai"; C = "https://coinranking.com/coin/MoTuySvg7+dai-dai"; D = "1.00"; E = "  +0.02%  " },
  @{ Row = 38; B = "InjectiveProtocol"; C = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; D = "37.19"; E = "  +2.21%  " },
  @{ Row = 39; B = "TheGraph"; C = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"; D = "0.379"; E = "  +0.12%  " },
  @{ Row = 40; B = "PEPE"; C = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"; D = "0.0₃0719"; E = "  +15.13%  " },
  @{ Row = 41; B = "FirstDigitalUSD"; C = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"; D = "0.996"; E = "  -0.10%  " },
  @{ Row = 42; B = "Kaspa"; C = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"; D = "0.124"; E = "  -0.12%  " },
  @{ Row = 43; B = "Maker"; C = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"; D = "2.887.48"; E = "  +3.22%  " },
  @{ Row = 44; B = "Fetch.AI"; C = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"; D = "2.54"; E = "  +9.13%  " },
  @{ Row = 45; B = "Stacks"; C = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"; D = "2.93"; E = "  +10.94%  " },
  @{ Row = 46; B = "WEMIXToken"; C = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"; D = "2.66"; E = "  +2.47%  " },
  @{ Row = 47; B = "VeChain"; C = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D = "0.0395"; E = "  +0.60%  " },
  @{ Row = 48; B = "ThetaToken"; C = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"; D = "2.61"; E = "  -4.80%  " },
  @{ Row = 49; B = "Stellar"; C = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D = "0.124"; E = "  +0.04%  " },
  @{ Row = 50; B = "Monero"; C = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D = "135.01"; E = "  +0.10%  " },
  @{ Row = 51; B = "ApeXProtocol"; C = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"; D = "2.92"; E = "  +6.30%  " }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("B$row").Value = "'" + $r.B
    $ws.Range("B$row").Style = "Normal"
    $ws.Range("C$row").Value = "'" + $r.C
    $ws.Range("C$row").Style = "Normal"
    $ws.Range("D$row").Value = "'" + $r.D
    $ws.Range("D$row").Style = "Normal"
    $ws.Range("E$row").Value = "'" + $r.E
    $ws.Range("E$row").Style = "Normal"
}

Write-Output "done"